$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: "Operasjonleie" -> "Avdekning/incisjon" + new helper text / logic columns ---
$ws.Cells.Item(13, 3).Value = "Avdekning/incisjon"
$ws.Cells.Item(13, 4).Value = "Legge inn link til tilsvarende dokument i DocMap (Ståle skulle underøke, håper det er avklart nå?)"

# --- "BNA" column split into separate, more precise labels ---
$ws.Cells.Item(2, 5).Value = "BNAMGR"
$ws.Cells.Item(12, 5).Value = "BNA/MGR"

# --- New content rows describing additional operator-assessment fields ---
$ws.Cells.Item(15, 3).Value = "Antibiotikaprofylakse"
$ws.Cells.Item(15, 4).Value = "Kan vi få en link til nasjonal veileder"
$ws.Cells.Item(16, 4).Value = "Kan vi få en link til nasjonal veileder"
$ws.Cells.Item(16, 3).Value = "Tromsboseprofylakse"

$ws.Cells.Item(14, 2).Value = "Spesille behov for kirurgi"
$ws.Cells.Item(14, 3).Value = "Alle elementer som er brukt"
$ws.Cells.Item(17, 3).Value = "Alle elementer som er brukt"
$ws.Cells.Item(18, 3).Value = "Alle elementer som er brukt"

$ws.Cells.Item(14, 4).Value = "Re-use. Gjenbruk informasjon fra beslutningsnotat."
$ws.Cells.Item(17, 4).Value = "Re-use. Gjenbruk informasjon fra beslutningsnotat."
$ws.Cells.Item(18, 4).Value = "Re-use. Gjenbruk informasjon fra beslutningsnotat."

$ws.Cells.Item(15, 2).Value = "Medikamenter ved krirugi"
$ws.Cells.Item(17, 2).Value = "Tilleggsinformasjon kirurgi"
$ws.Cells.Item(18, 2).Value = "Asa fysisk status klassifikasjon"

# --- Arketype column ("MGR") for the rows above ---
$ws.Cells.Item(13, 5).Value = "MGR"
$ws.Cells.Item(14, 5).Value = "MGR"
$ws.Cells.Item(15, 5).Value = "MGR"
$ws.Cells.Item(16, 5).Value = "MGR"
$ws.Cells.Item(17, 5).Value = "MGR"
$ws.Cells.Item(18, 5).Value = "MGR"

# Row 13 wraps onto two lines like rows 3/8/11
$ws.Rows.Item(13).RowHeight = 28.8

# Reflect where the author's cursor ended up after editing
$ws.Range("B25").Select()
